$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Fixed bug: grid mapping not deleting correct item - the authors/keywords
# pipe-delimited values for rows 2 and 3 were pointing at the wrong (blank)
# shared string; update them to the correct multi-value strings.
$ws1.Range("F2").Value = "a3|a1|a2"
$ws1.Range("G2").Value = "k3|k1|k2"
$ws1.Range("F3").Value = "a3|a1|a2"
$ws1.Range("G3").Value = "k3|k1|k2"

# Add Include all versions -> new "Sheet2" for the Error Messaging grid,
# inserted right after Sheet1.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

$headers = @("a_archive","a_compound_architecture","a_link_resolved","i_reference_cnt","i_has_folder","i_folder_id","r_composite_id","r_composite_label","r_component_label","r_order_no")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = [char](65 + $i)
    $cellRef = "${col}1"
    $cell = $ws2.Range($cellRef)
    $cell.Value = $headers[$i]
    $cell.NumberFormat = "@"
}

# Restore Sheet1's selection to the error-messaging header range and
# clear the old scroll/tab-selected state (now that Sheet2 is active).
$ws1.Activate()
$ws1.Range("Q1:Z1").Select()

# Sheet2 becomes the active/selected sheet with its header row selected.
$ws2.Activate()
$ws2.Range("A1:J1").Select()
